# Scheduled runner update: refresh cached market-board figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns)
# across the per-job Leve profit tables.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 106
$ws.Range("H106").Value = 2963.6667
$ws.Range("I106").Value = 2956.4
$ws.Range("K106").Value = 2956.4
$ws.Range("M106").Value = -2325.4
# Row 135
$ws.Range("H135").Value = 995.0714
$ws.Range("I135").Value = 748.7778
$ws.Range("K135").Value = 6739.000199999999
$ws.Range("M135").Value = -4204.000199999999
# Row 138
$ws.Range("H138").Value = 15629083
$ws.Range("I138").Value = 895.0714
$ws.Range("K138").Value = 2685.2142
$ws.Range("M138").Value = 2454.7858
# Row 141
$ws.Range("H141").Value = 1930
$ws.Range("I141").Value = 1342.5
$ws.Range("J141").Value = 3105
$ws.Range("K141").Value = 4027.5
$ws.Range("L141").Value = 9315
$ws.Range("M141").Value = 1152.5
$ws.Range("N141").Value = -19675

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3867.3125
$ws.Range("I45").Value = 3039.6667
$ws.Range("J45").Value = 4931.4287
$ws.Range("K45").Value = 3039.6667
$ws.Range("L45").Value = 4931.4287
$ws.Range("M45").Value = -2662.6667
$ws.Range("N45").Value = -5685.4287
# Row 61
$ws.Range("H61").Value = 2838.2173
$ws.Range("I61").Value = 2739.9546
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2739.9546
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2527.9546
$ws.Range("N61").Value = -5424
# Row 132
$ws.Range("H132").Value = 65176.285
$ws.Range("I132").Value = 8982.594999999999
$ws.Range("J132").Value = 481009.6
$ws.Range("K132").Value = 26947.785
$ws.Range("L132").Value = 1443028.8
$ws.Range("M132").Value = -24417.785
$ws.Range("N132").Value = -1448088.8
# Row 136
$ws.Range("H136").Value = 2838.2173
$ws.Range("I136").Value = 2739.9546
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8219.863799999999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -5669.863799999999
$ws.Range("N136").Value = -20100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1382.3636
$ws.Range("J94").Value = 2154.5
$ws.Range("L94").Value = 2154.5
$ws.Range("N94").Value = -3056.5
# Row 105
$ws.Range("H105").Value = 12214.871
$ws.Range("J105").Value = 6972.3687
$ws.Range("L105").Value = 6972.3687
$ws.Range("N105").Value = -10466.3687
# Row 134
$ws.Range("H134").Value = 1830
$ws.Range("I134").Value = 922.3333
$ws.Range("K134").Value = 2766.9999
$ws.Range("M134").Value = -231.9998999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
# Row 141
$ws.Range("H141").Value = 86441.14999999999
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 89148.58
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 89148.58
$ws.Range("M141").Value = -29820
$ws.Range("N141").Value = -99508.58

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 925
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 1066.6666
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 1066.6666
$ws.Range("M3").Value = -384
$ws.Range("N3").Value = -1298.6666
# Row 11
$ws.Range("H11").Value = 10482470
$ws.Range("I11").Value = 11876133
$ws.Range("J11").Value = 30000
$ws.Range("K11").Value = 11876133
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = -11875994
$ws.Range("N11").Value = -30278
# Row 14
$ws.Range("H14").Value = 2272181.8
$ws.Range("I14").Value = 3970794.5
$ws.Range("J14").Value = 7365
$ws.Range("K14").Value = 3970794.5
$ws.Range("L14").Value = 7365
$ws.Range("M14").Value = -3970626.5
$ws.Range("N14").Value = -7701
# Row 21
$ws.Range("H21").Value = 15242
$ws.Range("I21").Value = 15242
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 15242
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -15069
$ws.Range("N21").ClearContents()
# Row 30
$ws.Range("H30").Value = 15242
$ws.Range("I30").Value = 15242
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 15242
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -15137
$ws.Range("N30").ClearContents()
# Row 132
$ws.Range("H132").Value = 1648.2
$ws.Range("I132").Value = 1311.75
$ws.Range("K132").Value = 3935.25
$ws.Range("M132").Value = -1405.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3890.625
$ws.Range("I7").Value = 2244.8
$ws.Range("K7").Value = 2244.8
$ws.Range("M7").Value = -2132.8
# Row 23
$ws.Range("H23").Value = 6000
$ws.Range("I23").Value = 6000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -5770
$ws.Range("N23").ClearContents()
# Row 55
$ws.Range("H55").Value = 906.8214
$ws.Range("I55").Value = 611
$ws.Range("J55").Value = 1248.1538
$ws.Range("K55").Value = 611
$ws.Range("L55").Value = 1248.1538
$ws.Range("M55").Value = -438
$ws.Range("N55").Value = -1594.1538
# Row 122
$ws.Range("H122").Value = 3729.75
$ws.Range("I122").Value = 3397.7693
$ws.Range("K122").Value = 10193.3079
$ws.Range("M122").Value = -7743.3079
# Row 126
$ws.Range("H126").Value = 3890.625
$ws.Range("I126").Value = 2244.8
$ws.Range("K126").Value = 6734.400000000001
$ws.Range("M126").Value = -4264.400000000001
# Row 131
$ws.Range("H131").Value = 89993.336
$ws.Range("J131").Value = 89993.336
$ws.Range("L131").Value = 89993.336
$ws.Range("N131").Value = -100073.336
# Row 132
$ws.Range("H132").Value = 1958.3334
$ws.Range("I132").Value = 2022.909
$ws.Range("K132").Value = 6068.727000000001
$ws.Range("M132").Value = -3538.727000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 7617.864
$ws.Range("I14").Value = 118.6
$ws.Range("K14").Value = 118.6
$ws.Range("M14").Value = 49.40000000000001
# Row 18
$ws.Range("H18").Value = 137123.47
$ws.Range("I18").Value = 251658
$ws.Range("K18").Value = 251658
$ws.Range("M18").Value = -251485
# Row 126
$ws.Range("H126").Value = 16673643
$ws.Range("I126").Value = 23817920
$ws.Range("K126").Value = 71453760
$ws.Range("M126").Value = -71451290
# Row 132
$ws.Range("H132").Value = 2299.32
$ws.Range("I132").Value = 956.9167
$ws.Range("K132").Value = 2870.7501
$ws.Range("M132").Value = -340.7501000000002
# Row 136
$ws.Range("H136").Value = 3014.535
$ws.Range("I136").Value = 1863.7142
$ws.Range("J136").Value = 8049.375
$ws.Range("K136").Value = 5591.142599999999
$ws.Range("L136").Value = 24148.125
$ws.Range("M136").Value = -3041.142599999999
$ws.Range("N136").Value = -29248.125
